$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 08:47:40"
$ws.Range("H2").Value = "90%"
$ws.Range("K2").Value = "0.1 MJ/m2"
$ws.Range("E3").Value = "2026-02-07 08:47:42"
$ws.Range("K3").Value = "0.1 MJ/m2"
$ws.Range("N3").Value = "-7.9 °C 8:13 TU"
$ws.Range("E4").Value = "2026-02-07 08:47:45"
$ws.Range("J4").Value = "1001.9 hPa"
$ws.Range("K4").Value = "0.2 MJ/m2"
$ws.Range("O4").Value = "10.8 °C"
$ws.Range("E5").Value = "2026-02-07 08:47:47"
$ws.Range("J5").Value = "1001.9 hPa"
$ws.Range("K5").Value = "0.2 MJ/m2"
$ws.Range("O5").Value = "8.0 °C"
$ws.Range("E6").Value = "2026-02-07 08:47:50"
$ws.Range("H6").Value = "56%"
$ws.Range("J6").Value = "1003.4 hPa"
$ws.Range("K6").Value = "0.5 MJ/m2"
$ws.Range("E7").Value = "2026-02-07 08:47:52"
$ws.Range("J7").Value = "1003.1 hPa"
$ws.Range("K7").Value = "0.5 MJ/m2"
$ws.Range("E8").Value = "2026-02-07 08:47:54"
$ws.Range("H8").Value = "93%"
$ws.Range("K8").Value = "0.6 MJ/m2"
$ws.Range("M8").Value = "10.2 °C 8:14 TU"
$ws.Range("O8").Value = "4.3 °C"
$ws.Range("E9").Value = "2026-02-07 08:47:57"
$ws.Range("O9").Value = "1.1 °C"
$ws.Range("E10").Value = "2026-02-07 08:47:59"
$ws.Range("H10").Value = "97%"
$ws.Range("M10").Value = "11.9 °C 8:29 TU"
$ws.Range("O10").Value = "7.6 °C"
$ws.Range("E11").Value = "2026-02-07 08:48:01"
$ws.Range("J11").Value = "1005.9 hPa"
$ws.Range("K11").Value = "0.1 MJ/m2"
$ws.Range("E12").Value = "2026-02-07 08:48:04"
$ws.Range("H12").Value = "67%"
$ws.Range("K12").Value = "0.4 MJ/m2"
$ws.Range("M12").Value = "12.6 °C 8:27 TU"
$ws.Range("O12").Value = "10.0 °C"
$ws.Range("E13").Value = "2026-02-07 08:48:06"
$ws.Range("H13").Value = "82%"
$ws.Range("M13").Value = "13.2 °C 8:25 TU"
$ws.Range("O13").Value = "8.4 °C"
$ws.Range("E14").Value = "2026-02-07 08:48:08"
$ws.Range("I14").Value = "0.3 mm"
$ws.Range("O14").Value = "-6.1 °C"
$ws.Range("E15").Value = "2026-02-07 08:48:11"
$ws.Range("J15").Value = "1002.2 hPa"
$ws.Range("K15").Value = "0.3 MJ/m2"
$ws.Range("E16").Value = "2026-02-07 08:48:13"
$ws.Range("K16").Value = "0.1 MJ/m2"
$ws.Range("E17").Value = "2026-02-07 08:48:16"
$ws.Range("J17").Value = "1005.4 hPa"
$ws.Range("K17").Value = "0.1 MJ/m2"
$ws.Range("O17").Value = "3.0 °C"
$ws.Range("E18").Value = "2026-02-07 08:48:18"
$ws.Range("K18").Value = "0.3 MJ/m2"
$ws.Range("O18").Value = "-7.8 °C"
$ws.Range("E19").Value = "2026-02-07 08:48:21"
$ws.Range("J19").Value = "1006.8 hPa"
$ws.Range("K19").Value = "0.5 MJ/m2"
$ws.Range("M19").Value = "5.8 °C 8:29 TU"
$ws.Range("O19").Value = "3.8 °C"
$ws.Range("E20").Value = "2026-02-07 08:48:23"
$ws.Range("H20").Value = "79%"
$ws.Range("K20").Value = "0.3 MJ/m2"
$ws.Range("M20").Value = "-3.4 °C 8:11 TU"
$ws.Range("E21").Value = "2026-02-07 08:48:26"
$ws.Range("J21").Value = "1002.6 hPa"
$ws.Range("K21").Value = "0.3 MJ/m2"
$ws.Range("E22").Value = "2026-02-07 08:48:28"
$ws.Range("H22").Value = "87%"
$ws.Range("K22").Value = "0.7 MJ/m2"
$ws.Range("M22").Value = "12.6 °C 8:29 TU"
$ws.Range("O22").Value = "7.1 °C"
$ws.Range("E23").Value = "2026-02-07 08:48:31"
$ws.Range("H23").Value = "95%"
$ws.Range("J23").Value = "1002.0 hPa"
$ws.Range("K23").Value = "0.5 MJ/m2"
$ws.Range("M23").Value = "10.6 °C 8:29 TU"
$ws.Range("O23").Value = "7.6 °C"
$ws.Range("E24").Value = "2026-02-07 08:48:33"
$ws.Range("H24").Value = "79%"
$ws.Range("J24").Value = "1001.2 hPa"
$ws.Range("K24").Value = "0.5 MJ/m2"
$ws.Range("E25").Value = "2026-02-07 08:48:35"
$ws.Range("H25").Value = "96%"
$ws.Range("J25").Value = "1005.6 hPa"
$ws.Range("K25").Value = "0.1 MJ/m2"
$ws.Range("E26").Value = "2026-02-07 08:48:38"
$ws.Range("O26").Value = "-2.7 °C"
$ws.Range("E27").Value = "2026-02-07 08:48:40"
$ws.Range("H27").Value = "89%"
$ws.Range("J27").Value = "1001.8 hPa"
$ws.Range("K27").Value = "0.4 MJ/m2"
$ws.Range("L27").Value = "37.8 km/h - 346º 8:16 TU"
$ws.Range("M27").Value = "12.4 °C 8:28 TU"
$ws.Range("O27").Value = "9.4 °C"
$ws.Range("E28").Value = "2026-02-07 08:48:43"
$ws.Range("J28").Value = "1004.7 hPa"
$ws.Range("O28").Value = "2.7 °C"
$ws.Range("E29").Value = "2026-02-07 08:48:45"
$ws.Range("H29").Value = "60%"
$ws.Range("K29").Value = "0.4 MJ/m2"
$ws.Range("O29").Value = "10.5 °C"
$ws.Range("E30").Value = "2026-02-07 08:48:47"
$ws.Range("H30").Value = "77%"
$ws.Range("K30").Value = "0.5 MJ/m2"
$ws.Range("E31").Value = "2026-02-07 08:48:50"
$ws.Range("J31").Value = "1006.3 hPa"
$ws.Range("E32").Value = "2026-02-07 08:48:52"
$ws.Range("J32").Value = "1004.9 hPa"
$ws.Range("K32").Value = "0.4 MJ/m2"
$ws.Range("E33").Value = "2026-02-07 08:48:55"
$ws.Range("H33").Value = "90%"
$ws.Range("O33").Value = "6.9 °C"
$ws.Range("E34").Value = "2026-02-07 08:48:57"
$ws.Range("K34").Value = "0.2 MJ/m2"
$ws.Range("O34").Value = "5.9 °C"
$ws.Range("E35").Value = "2026-02-07 08:48:59"
$ws.Range("K35").Value = "0.1 MJ/m2"
$ws.Range("N35").Value = "-8.5 °C 8:27 TU"
$ws.Range("O35").Value = "-6.3 °C"
$ws.Range("E36").Value = "2026-02-07 08:49:02"
$ws.Range("J36").Value = "1007.2 hPa"
$ws.Range("K36").Value = "0.4 MJ/m2"
